$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")

# Rename existing "NormalUser" entry (row 3) to "UserOne"
$ws.Range("A3").Value = "UserOne"
$ws.Range("B3").Value = "UserOne"

# New users to append, following the same A=B=name, C=2040 pattern
$newUsers = @("UserTwo", "UserThree", "UserFour", "UserFive", "UserSix")

$row = 4
foreach ($name in $newUsers) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 3).Value = 2040
    $row++
}

# Match the selection state from the diff (activeCell B5)
$ws.Range("B5").Select()
